$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert two new columns before column B, shifting the existing
# Variable_Name/Question_Type/Required/List_Values/If_Condition/
# Then_Question/Else_Question columns from B..H to D..J.
$ws1.Range("B:C").EntireColumn.Insert() | Out-Null

# New column C = "Description" (filled first, matching entry order)
$ws1.Cells.Item(1,3).Value = "Description"

# New column B = "Title"
$ws1.Cells.Item(2,2).Value = "Patient"
$ws1.Cells.Item(3,2).Value = "Age"
$ws1.Cells.Item(4,2).Value = "Sex"
$ws1.Cells.Item(5,2).Value = "Pregnant"
$ws1.Cells.Item(6,2).Value = "symptoms"

$ws1.Cells.Item(3,3).Value = "Please enter Age …."
$ws1.Cells.Item(2,3).Value = "Please enter Name…."
$ws1.Cells.Item(4,3).Value = "Please enter Sex…."
$ws1.Cells.Item(5,3).Value = "Please enter Pregnant..."
$ws1.Cells.Item(6,3).Value = "Please enter symptoms…"

$ws1.Cells.Item(1,2).Value = "Title"

# Rename the shifted If_Condition follow-up headers.
$ws1.Cells.Item(1,9).Value = "Then_Goto"
$ws1.Cells.Item(1,10).Value = "Else_Goto"

# The boolean (Required) list validation sqref shifts from D to F
# automatically as part of the column insert above - nothing to do.

# Re-point the DataTypes list validation at its new column (E).
$ws1.Range("E2:E1048576").Validation.Add(3, 1, 1, "=DataTypes!`$A:`$A") | Out-Null

# Match the recorded selection on Sheet1.
$ws1.Range("I1").Select() | Out-Null
